$d = $word.ActiveDocument

# Determine the range to replace: from start of paragraph 5 ("News articles...")
# through end of the document content (paragraph 15, "I don't know what to add here.")
$p5 = $d.Paragraphs(5)
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$rng = $d.Range($p5.Range.Start, $lastPara.Range.End)

$bodyXml = '<w:p><w:r><w:t xml:space="preserve">We will be using the dataset from the Di-Tech challenge consisting of data collected around the ride-hailing company </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Didi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>Chuxing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">.  The dataset includes data regarding the transportation </w:t></w:r><w:r><w:t xml:space="preserve">behaviors Chinese citizens. </w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Structure:</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p><w:p><w:pPr><w:tabs><w:tab w:val="right" w:pos="9360"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">We intend on applying regression techniques to the dataset. </w:t></w:r><w:r><w:t xml:space="preserve"> Through applying these techniques, we will be able to mine traffic patterns of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Didi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Chuxing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> customers.  The traffic patterns</w:t></w:r><w:r><w:t xml:space="preserve"> we intend to data mine</w:t></w:r><w:r><w:t xml:space="preserve"> will include </w:t></w:r><w:r><w:t xml:space="preserve">data regarding the time of the day, location of pick-ups and drop-offs, and other relevant information to forecast </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Didi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Chuxing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> customer’s needs. </w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="right" w:pos="9360"/></w:tabs><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Problem:</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="right" w:pos="9360"/></w:tabs></w:pPr><w:r><w:t>With such a high reliance on non-personal transportation, it is important for the transportation industry to have a good understanding of the transportation patterns.  Using this data</w:t></w:r><w:r><w:t xml:space="preserve">set, we hope to forecast the transportation needs of the customers of </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Didi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Chuxing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">.  By focusing better on the needs of the customers, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Didi</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Chuxing</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> can optimize their drivers and their revenue.</w:t></w:r></w:p><w:p><w:pPr><w:tabs><w:tab w:val="right" w:pos="9360"/></w:tabs><w:rPr><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">What is </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:b/></w:rPr><w:t>New:</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:tabs><w:tab w:val="right" w:pos="9360"/></w:tabs></w:pPr><w:r><w:t xml:space="preserve">We will explore a deeper understanding of regression techniques by comparing several regression techniques to determine the advantages and disadvantages to the techniques.  </w:t></w:r><w:r><w:rPr><w:highlight w:val="yellow"/></w:rPr><w:t>(I need help here)</w:t></w:r></w:p>
'

$xml = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($xml)

Write-Output $d.Content.Text
